$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
for ($i=1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        Write-Host "$i : $($sh.Name) : [$($sh.TextFrame.TextRange.Text)]"
    }
}
